$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the two new header cells, copying the format of the existing H1 header
# so they share the same cell style (bold, bordered, centered) as the rest
# of the header row.
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)

# Fill in the data for the new I and J columns (rows 2-8).
$valuesI = @(9, 9, 9, 9, 7, 8, 8)
$valuesJ = @(9, 9, 9, 9, 7, 9, 8)

for ($i = 0; $i -lt 7; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 9).Value = $valuesI[$i]
    $ws.Cells.Item($row, 10).Value = $valuesJ[$i]
}
